$d = $word.ActiveDocument

function Insert-LineBreak($findText, $replaceText) {
    $found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Find.Execute could not find: $findText"
    }
}

# Paragraph: Objetivos
Insert-LineBreak "materiais.Desenvolver a competência" "materiais.^lDesenvolver a competência"
Insert-LineBreak "resultados analíticos. Incentivar trabalhos" "resultados analíticos. ^lIncentivar trabalhos"

# Paragraph: Programa resumido
Insert-LineBreak "Conformação Mecânica1.1. Classificação" "Conformação Mecânica^l1.1. Classificação"
Insert-LineBreak "Calandragem.2. Processos de Usinagem " "Calandragem.^l2. Processos de Usinagem "
Insert-LineBreak "2. Processos de Usinagem 2.1 Principais" "2. Processos de Usinagem ^l2.1 Principais"

# Paragraph: Programa
Insert-LineBreak "molecular.Espectroscopia UV/Visível" "molecular.^lEspectroscopia UV/Visível"

# Paragraph: Bibliografia
Insert-LineBreak "2003.2. SKOOG" "2003.^l2. SKOOG"
Insert-LineBreak "2007.3. MITRA" "2007.^l3. MITRA"
Insert-LineBreak "2003.4. ANDERSON" "2003.^l4. ANDERSON"
